# Append two new rows to Sheet3 with Name="Abiram", Country="Swouvania", Company="Uthesh"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

$ws.Cells.Item(62, 1).Value = "Abiram"
$ws.Cells.Item(62, 2).Value = "Swouvania"
$ws.Cells.Item(62, 3).Value = "Uthesh"

$ws.Cells.Item(63, 1).Value = "Abiram"
$ws.Cells.Item(63, 2).Value = "Swouvania"
$ws.Cells.Item(63, 3).Value = "Uthesh"
